$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Actual Consumption (MW)" values for rows 2..34 (column A)
$aValues = @(6412,6341,6293,6234,6261,6206,6209,6166,6104,6065,6056,6038,5994,5969,5990,5980,5994,5976,5992,6006,6040,6061,6095,6130,6167,6202,6241,6294,6289,6472,6551,6735,6884)

# New "Timestamp" serial values for rows 2..34 (column B)
$bValues = @(46073.95833333334,46073.96875,46073.97916666666,46073.98958333334,46074,46074.01041666666,46074.02083333334,46074.03125,46074.04166666666,46074.05208333334,46074.0625,46074.07291666666,46074.08333333334,46074.09375,46074.10416666666,46074.11458333334,46074.125,46074.13541666666,46074.14583333334,46074.15625,46074.16666666666,46074.17708333334,46074.1875,46074.19791666666,46074.20833333334,46074.21875,46074.22916666666,46074.23958333334,46074.25,46074.26041666666,46074.27083333334,46074.28125,46074.29166666666)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# Row 34 is newly added; copy the date number format style from row 33's B cell
$ws.Range("B33").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$ws.Cells.Item(34, 2).Value = $bValues[$aValues.Length - 1]
